# "client & group Scenarios"
# Remove the two rows that hold the standalone "charges"/"addcharges"
# add-charge helper entries (A30:B31) from the ProductLoanInput sheet,
# shifting the remaining rows (fundsource .. overpaymentliability) up by two.
# Excel prunes the now-unused shared strings ("charges", "addcharges",
# "SpecifiedDueDateFees-Flat") automatically on save, which renumbers the
# shared-string indices referenced elsewhere (e.g. the product-name header
# in B1 of both sheets, and the allowPartialPeriodInterestCalcualtion label).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 30 and 31 (charges/SpecifiedDueDateFees-Flat, addcharges/checked);
# everything below shifts up automatically.
$ws.Rows("30:31").Delete() | Out-Null

# Match the author's recorded selection/scroll state after the edit.
$ws.Range("A33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
